$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the touched columns so numeric-looking
# strings (prices, percentages, hour values) are stored as literal text,
# matching the source data which uses inline string cells.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "309.55"
$ws.Range("E2").Value = "-2.66%"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "52.22"
$ws.Range("E3").Value = "9.08%"
$ws.Range("G3").Value = "14"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").Value = "-2.32%"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.07802"
$ws.Range("E5").Value = "-1.97%"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "4.505"
$ws.Range("E6").Value = "-2.10%"
$ws.Range("G6").Value = "14"
$ws.Range("D7").Value = "1.365"
$ws.Range("E7").Value = "-4.02%"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "1.582"
$ws.Range("E8").Value = "-3.59%"
$ws.Range("G8").Value = "14"
$ws.Range("D9").Value = "0.1227"
$ws.Range("E9").Value = "-3.78%"
$ws.Range("G9").Value = "14"
$ws.Range("D10").Value = "0.2008"
$ws.Range("E10").Value = "3.81%"
$ws.Range("G10").Value = "14"
$ws.Range("D11").Value = "0.09671"
$ws.Range("E11").Value = "4.34%"
$ws.Range("G11").Value = "14"
$ws.Range("D12").Value = "0.04723"
$ws.Range("E12").Value = "3.13%"
$ws.Range("G12").Value = "14"
$ws.Range("E13").Value = "0.14%"
$ws.Range("G13").Value = "14"
$ws.Range("D14").Value = "0.001282"
$ws.Range("E14").Value = "-4.66%"
$ws.Range("G14").Value = "14"
$ws.Range("D15").Value = "0.005789"
$ws.Range("E15").Value = "-1.31%"
$ws.Range("G15").Value = "14"
$ws.Range("E16").Value = "2,006.97%"
$ws.Range("G16").Value = "14"
$ws.Range("D17").Value = "3.335"
$ws.Range("E17").Value = "0.22%"
$ws.Range("G17").Value = "14"
$ws.Range("E18").Value = "-0.83%"
$ws.Range("G18").Value = "14"
$ws.Range("D19").Value = "0.3443"
$ws.Range("E19").Value = "-0.68%"
$ws.Range("G19").Value = "14"
$ws.Range("D20").Value = "7.998"
$ws.Range("E20").Value = "-1.42%"
$ws.Range("G20").Value = "14"
$ws.Range("D21").Value = "0.1369"
$ws.Range("E21").Value = "-1.99%"
$ws.Range("G21").Value = "14"
$ws.Range("D22").Value = "0.3093"
$ws.Range("E22").Value = "-0.32%"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "0.04168"
$ws.Range("E23").Value = "-0.19%"
$ws.Range("G23").Value = "14"
$ws.Range("E24").Value = "-3.94%"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.003950"
$ws.Range("E25").Value = "-6.75%"
$ws.Range("G25").Value = "14"
$ws.Range("D26").Value = "0.0001351"
$ws.Range("E26").Value = "-0.01%"
$ws.Range("G26").Value = "14"
$ws.Range("G27").Value = "14"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("D38").Value = "0.02601"
$ws.Range("E38").Value = "-1.06%"
$ws.Range("G38").Value = "14"
$ws.Range("D39").Value = "0.05890"
$ws.Range("E39").Value = "4.32%"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.01100"
$ws.Range("E40").Value = "4.82%"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.007900"
$ws.Range("E41").Value = "-1.52%"
$ws.Range("G41").Value = "14"
$ws.Range("D42").Value = "0.1421"
$ws.Range("E42").Value = "-0.84%"
$ws.Range("G42").Value = "14"
$ws.Range("D43").Value = "0.008238"
$ws.Range("E43").Value = "6.91%"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.008454"
$ws.Range("E44").Value = "-0.60%"
$ws.Range("G44").Value = "14"
$ws.Range("D45").Value = "0.3113"
$ws.Range("E45").Value = "-10.17%"
$ws.Range("G45").Value = "14"
$ws.Range("D46").Value = "0.00007355"
$ws.Range("E46").Value = "6.20%"
$ws.Range("G46").Value = "14"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.33%"
$ws.Range("G47").Value = "14"
$ws.Range("D48").Value = "0.05692"
$ws.Range("E48").Value = "3.83%"
$ws.Range("G48").Value = "14"
$ws.Range("D49").Value = "0.002621"
$ws.Range("E49").Value = "-34.73%"
$ws.Range("G49").Value = "14"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.33%"
$ws.Range("G50").Value = "14"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.33%"
$ws.Range("G51").Value = "14"

# Restore the default (General) style on the touched columns so the
# cells keep their original look-and-feel; only the values changed.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
